# Transportation Technology Shareweights.xlsx - update sales-share input
# cells on the "Data" sheet (rows 24-30). A handful of these cells were
# formulas pulling from 'SYVbT-passenger' / other Data cells; the author
# replaced them with plain hard-coded numeric values (and changed a few
# plain numeric inputs too), to smooth the sales shares used by the
# S-curve / linear interpolation further to the right on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 24: plain numeric inputs
$ws.Range("D24").Value = 0.6
$ws.Range("E24").Value = 0.6
$ws.Range("F24").Value = 1

# Row 25: both E25 and F25 were formulas; replace with literal values
$ws.Range("E25").Value = 0.29802375741500553
$ws.Range("F25").Value = 0.15850139443373243

# Row 26: E26 was a plain value, F26 was "=E26"; both become literals
$ws.Range("E26").Value = 0.29802375741500553
$ws.Range("F26").Value = 0.29802375741500553

# Row 27: plain numeric inputs
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 1

# Row 28: plain numeric inputs
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0.006

# Row 29: E29 was a formula; F29 was a plain value
$ws.Range("E29").Value = 0.021496445375763083
$ws.Range("F29").Value = 0.53425472287359943

# Row 30: E30 was a formula; F30 is unchanged (0.05) so left untouched
$ws.Range("E30").Value = 0.00008470448323552864

# Reflect the author's final cursor position / scroll on the Data sheet
$ws.Activate()
$ws.Range("C26").Select()
$excel.ActiveWindow.ScrollRow = 15
